$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted right before the current row 686,
# pushing all subsequent rows (old 686-702) down by one (new 687-703).
$ws.Rows.Item(686).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A686").Value = 3
$ws.Range("B686").Value = "Femacal de La Calera"
$ws.Range("C686").Value = "Coquimbo"
$ws.Range("D686").Value = 45239
$ws.Range("E686").Value = 5
$ws.Range("F686").Value = 100112031
$ws.Range("G686").Value = "Poroto verde"
$ws.Range("H686").Value = "Sin especificar"
$ws.Range("I686").Value = "Primera"
$ws.Range("J686").Value = 50
$ws.Range("K686").Value = 54000
$ws.Range("L686").Value = 55000
$ws.Range("M686").Value = 54500
$ws.Range("N686").Value = "`$/malla 25 kilos"
$ws.Range("O686").Value = "Provincia de Limarí"
$ws.Range("P686").Value = 2180
$ws.Range("Q686").Value = 25
$ws.Range("R686").Value = "Hortaliza"
